$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 204; $r -le 303; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "intellectual property rights") {
        $cell.Value2 = "intellectual property"
    }
}
